# edit.ps1
# Applies the changes described by the target diff to the PSP time-recording
# worksheet: fills in the new "10/17" log entry in row 21, right-aligns the
# date cell A19 (new cell style), updates the saved cell selection to A19,
# and adds the two new shared strings (with their original mixed Korean /
# English run-level font formatting reproduced as closely as possible).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper functions to apply the "Korean UI font" (돋움) and the "Latin/number
# font" (Arial) to a run of characters within a cell's text, matching the
# rich-text runs originally produced by Excel for this workbook.
function Set-KoreanRun($range, $start, $len) {
    $f = $range.Characters($start, $len).Font
    $f.Name = "돋움"
    $f.Size = 10
    $f.ColorIndex = -4105
}
function Set-LatinRun($range, $start, $len) {
    $f = $range.Characters($start, $len).Font
    $f.Name = "Arial"
    $f.Size = 10
    $f.ColorIndex = -4105
}

# ---------------------------------------------------------------------------
# Row 21: new time-log entry for "10월 17일" (Oct 17), 17:00-21:00, 43 min
# interruption, 240 min delta, for "Initial Data Set 자료 수집을 위한 설문지
# 제작" (Preparing survey for Initial Data Set collection).
# ---------------------------------------------------------------------------

$a21 = $ws.Range("A21")
$a21.Value = "10월 17일"
Set-KoreanRun $a21 3 1
Set-LatinRun  $a21 4 3
Set-KoreanRun $a21 7 1

$ws.Range("B21").Value = 0.70833333333333337
$ws.Range("C21").Value = 0.875
$ws.Range("D21").Value = 43
$ws.Range("E21").Value = 240

$f21 = $ws.Range("F21")
$f21.Value = "Initial Data Set 자료 수집을 위한 설문지 제작"
Set-KoreanRun $f21 18 2
Set-LatinRun  $f21 20 1
Set-KoreanRun $f21 21 3
Set-LatinRun  $f21 24 1
Set-KoreanRun $f21 25 2
Set-LatinRun  $f21 27 1
Set-KoreanRun $f21 28 3
Set-LatinRun  $f21 31 1
Set-KoreanRun $f21 32 2

# ---------------------------------------------------------------------------
# A19: right-align the date text (creates the new cell style used by the
# workbook after the edit).
# ---------------------------------------------------------------------------
$ws.Range("A19").HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# Update the sheet's saved selection to A19 (was F21).
# ---------------------------------------------------------------------------
$ws.Range("A19").Select()
